$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.491.77"
$ws.Range("E2").Value = "  -1.00%  "
$ws.Range("D3").Value = "3.971.98"
$ws.Range("E3").Value = "  -2.58%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "539.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.09%  "
$ws.Range("D7").Value = "3.963.85"
$ws.Range("E7").Value = "  -2.58%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.690"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("E10").Value = "  -3.82%  "
$ws.Range("E11").Value = "  -6.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +13.14%  "
$ws.Range("E13").Value = "  -4.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.75%  "
$ws.Range("D15").Value = "4.600.32"
$ws.Range("E15").Value = "  -2.70%  "
$ws.Range("D16").Value = "3.970.04"
$ws.Range("E16").Value = "  -2.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.07"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.55%  "
$ws.Range("E19").Value = "  -1.71%  "
$ws.Range("E20").Value = "  -5.32%  "
$ws.Range("D21").Value = "71.328.47"
$ws.Range("E21").Value = "  -1.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "431.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.75%  "
$ws.Range("E23").Value = "  -0.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "97.54"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.47"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +20.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.89%  "
$ws.Range("E30").Value = "  +1.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +19.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "51.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +20.74%  "
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "682.93"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "65.67"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.444"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.13%  "
$ws.Range("D39").Value = "0.0₃0830"
$ws.Range("E39").Value = "  -7.06%  "
$ws.Range("E40").Value = "  -3.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.77%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("E44").Value = "  -3.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.23"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.39"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.150"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.00%  "
$ws.Range("E48").Value = "  -1.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.03"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000270"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.20%  "
